$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove "DTD" suffix from player names (day-to-day injury tag no longer tracked).
# Excel appends newly-typed unique strings to the shared-string table, so setting
# these cell values in ascending row order reproduces the same shared-string layout
# as the authored edit (stale "...DTD" entries become orphaned, clean names land at
# the end of the table in this exact order).
$ws.Range("A2").Value2 = 'Giannis Antetokounmpo'
$ws.Range("A8").Value2 = 'Russell Westbrook'
$ws.Range("A18").Value2 = 'Jimmy Butler'
$ws.Range("A32").Value2 = 'Draymond Green'
$ws.Range("A35").Value2 = 'Devin Booker'
$ws.Range("A36").Value2 = 'Lou Williams'
$ws.Range("A37").Value2 = 'Clint Capela'
$ws.Range("A40").Value2 = 'Dwight Howard'
$ws.Range("A42").Value2 = 'Julius Randle'
$ws.Range("A64").Value2 = 'Jamal Murray'
$ws.Range("A69").Value2 = 'Harrison Barnes'
$ws.Range("A90").Value2 = 'Reggie Jackson'
$ws.Range("A92").Value2 = 'Gary Harris'
$ws.Range("A94").Value2 = 'Rondae Hollis-Jefferson'
$ws.Range("A104").Value2 = 'James Johnson'
$ws.Range("A107").Value2 = 'Nikola Mirotic'
$ws.Range("A118").Value2 = 'Dwight Powell'
$ws.Range("A127").Value2 = 'Josh Richardson'
$ws.Range("A149").Value2 = 'J.J. Barea'
$ws.Range("A159").Value2 = 'Thaddeus Young'
$ws.Range("A174").Value2 = 'Jonathan Isaac'
$ws.Range("A187").Value2 = 'Michael Kidd-Gilchrist'
$ws.Range("A188").Value2 = 'Dewayne Dedmon'
$ws.Range("A190").Value2 = 'Bam Adebayo'
$ws.Range("A202").Value2 = 'Courtney Lee'
$ws.Range("A217").Value2 = 'Wayne Ellington'
$ws.Range("A219").Value2 = 'Milos Teodosic'
$ws.Range("A220").Value2 = 'Maurice Harkless'
$ws.Range("A235").Value2 = 'Nemanja Bjelica'
$ws.Range("A255").Value2 = 'Lonnie Walker IV'
$ws.Range("A259").Value2 = 'PJ Tucker'
$ws.Range("A273").Value2 = 'Denzel Valentine'
$ws.Range("A290").Value2 = 'MarShon Brooks'
$ws.Range("A321").Value2 = 'Zhou Qi'
$ws.Range("A328").Value2 = 'Shabazz Napier'
$ws.Range("A329").Value2 = 'Omari Spellman'
$ws.Range("A334").Value2 = 'Jon Leuer'
$ws.Range("A348").Value2 = 'Frank Jackson'
$ws.Range("A351").Value2 = 'Terrance Ferguson'
$ws.Range("A369").Value2 = 'Iman Shumpert'
$ws.Range("A373").Value2 = 'Darius Miller'
$ws.Range("A376").Value2 = 'Patrick McCaw'
$ws.Range("A403").Value2 = 'Sterling Brown'
$ws.Range("A426").Value2 = 'Shake Milton'
$ws.Range("A430").Value2 = 'Alan Williams'
$ws.Range("A436").Value2 = 'Nene Hilario'
$ws.Range("A439").Value2 = 'Alex Abrines'
$ws.Range("A449").Value2 = 'Isaiah Canaan'
$ws.Range("A476").Value2 = 'Omer Asik'

# Restore the view/selection state recorded in the saved workbook.
$ws.Range("A13").Select()
$win = $excel.ActiveWindow
$win.ScrollRow = 17
$win.ScrollColumn = 1
